$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(22,'s21','s21_IMG_3174.jpeg','meltpatch','1479','525','104','52','173','2'),
    @(23,'s22','s22_IMG_3179.jpeg','meltpatch','328','1975','104','52','21','2'),
    @(24,'s23','s23_IMG_3176.jpeg','meltpatch','148','1863','104','52','169','2'),
    @(25,'s24','s24_IMG_3178.jpeg','meltpatch','2952','2329','104','52','78','2'),
    @(26,'s25','s25_IMG_3175.jpeg','meltpatch','851','822','104','52','136','2'),
    @(27,'s26','s26_IMG_3179.jpeg','meltpatch','833','123','104','52','142','2'),
    @(28,'s27','s27_IMG_3175.jpeg','meltpatch','2136','1530','104','52','30','2'),
    @(29,'s28','s28_IMG_3180.jpeg','meltpatch','2666','616','104','52','65','2'),
    @(30,'s29','s29_IMG_3178.jpeg','meltpatch','2137','1972','104','52','172','2'),
    @(31,'s30','s30_IMG_3174.jpeg','meltpatch','273','2264','104','52','135','2'),
    @(32,'s31','s31_e43_321_1_1.jpeg','meltpatch','1544','420','104','52','171','2'),
    @(33,'s32','s32_e45_321_1_3.jpeg','meltpatch','2541','1628','104','52','10','2'),
    @(34,'s33','s33_e46_321_2_2.jpeg','meltpatch','805','651','104','52','58','2'),
    @(35,'s34','s34_e48_321_2_0.jpeg','meltpatch','1450','2250','104','52','70','2'),
    @(36,'s35','s35_e44_321_1_2.jpeg','meltpatch','183','2513','104','52','173','2'),
    @(37,'s36','s36_e54_321_2_1.jpeg','meltpatch','110','879','77','38','117','2'),
    @(38,'s37','s37_e52_321_1_3.jpeg','meltpatch','2005','1343','77','38','108','2'),
    @(39,'s38','s38_e53_321_2_2.jpeg','meltpatch','1750','1919','77','38','168','2'),
    @(40,'s39','s39_e56_321_3_1.jpeg','meltpatch','2280','1956','77','38','75','2'),
    @(41,'s40','s40_e50_321_1_1.jpeg','meltpatch','609','57','77','38','83','2'),
    @(42,'s41','s41_e59_321_1_3.jpeg','meltpatch','2465','1193','77','38','32','2'),
    @(43,'s42','s42_e57_321_1_1.jpeg','meltpatch','2608','577','77','38','144','2'),
    @(44,'s43','s43_e63_321_3_1.jpeg','meltpatch','2366','1022','77','38','134','2'),
    @(45,'s44','s44_e58_321_1_2.jpeg','meltpatch','2078','487','77','38','107','2'),
    @(46,'s45','s45_e60_321_2_2.jpeg','meltpatch','1798','1861','77','38','72','2')
)

$columns = @("A","B","C","D","E","F","G","H","I")

foreach ($rowEntry in $rowsData) {
    $rowNum = $rowEntry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $colLetter = $columns[$i]
        $cellValue = $rowEntry[$i + 1]
        $cellRef = "$colLetter$rowNum"
        if ($i -ge 3) {
            # Columns D-I hold numeric-looking text (ids/coords/tolerances/theta/ratio);
            # force text format first so Excel doesn't coerce them to numbers.
            $ws.Range($cellRef).NumberFormat = "@"
        }
        $ws.Range($cellRef).Value = $cellValue
    }
}
